# Updates the crypto price/volume table to the latest scraped snapshot.
# Mirrors the GitHub Actions "Updated cryptos list" commit: per-row Price (D)
# and Volume(1h) (E) refreshes, plus a few coins that changed rank (full row swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.622.08"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "1.798.17"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'227.24"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'32.89"
$ws.Range("E8").Value = "  +3.37%  "
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "2.056.91"
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.10"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.792.28"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("D16").Value = "34.589.14"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "'4.31"
$ws.Range("E17").Value = "  +3.06%  "
$ws.Range("D18").Value = "'68.98"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0804"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'247.59"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "'11.32"
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "'168.54"
$ws.Range("E24").Value = "  +3.57%  "
$ws.Range("D25").Value = "'2.06"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'4.09"
$ws.Range("E30").Value = "  +11.28%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.24"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.0525"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("E33").Value = "  +1.79%  "
$ws.Range("E34").Value = "  +2.82%  "
$ws.Range("D35").Value = "1.428.78"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("E36").Value = "  +6.91%  "
$ws.Range("E37").Value = "  +2.65%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.06"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0193"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").Value = "'85.69"
$ws.Range("E40").Value = "  +6.74%  "
$ws.Range("D41").Value = "'2.41"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "'0.937"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("E43").Value = "  +3.24%  "
$ws.Range("D44").Value = "'13.78"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "1.956.77"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").Value = "'106.11"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("E51").Value = "  -5.66%  "
